$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Density" column header in K1
$ws.Range("K1").Value = "Density"

# Add density formulas (Post-Mix Mass / Post-Mix Vol) for the two populated data rows
$ws.Range("K3").Formula = "=I3/H3"
$ws.Range("K4").Formula = "=I4/H4"

# Update the active selection to K5, matching the new working cell
$ws.Range("K5").Select()
